$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column D header and data, factoring in the new "expand_copol" row (Val/Dev/Status)
$ws.Range("D6").Value = "Val"
$ws.Range("D7").Value = "Dev"
$ws.Range("D1").Value = "Status"

# Update the active selection to reflect where the user ended up after editing
$ws.Range("D1").Select()
